$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "'"
$ws.Range("B7").Value = "يامن "
$ws.Range("C7").Value = "'22"
$ws.Range("D7").Value = "الصمود"
$ws.Range("E7").Value = "الرحلة 3"
$ws.Range("F7").Value = "C3"
$ws.Range("G7").Value = "WCK"
$ws.Range("H7").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٤٨:٠١ م"
